$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "64.506.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -0.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.361.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -2.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "556.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -2.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "175.71"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  +0.67%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  -0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.354.61"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -2.10%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  +0.00%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  +3.73%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.630"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  +1.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "54.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  -0.78%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +0.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "9.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -0.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.889.15"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -2.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "18.43"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -1.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.360.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -2.06%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  +0.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "64.415.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -0.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.986"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -0.27%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "462.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  +13.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "4.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +11.93%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -2.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "86.12"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  +3.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "13.41"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +1.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  +1.44%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +2.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "8.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -1.69%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "30.12"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  +0.90%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  +0.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "11.48"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "580.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -1.74%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  +0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "58.89"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  -1.01%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  +0.15%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -8.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "3.50"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -0.62%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value2 = "PEPE"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0₃0758"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -1.25%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value2 = "InjectiveProtocol"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "35.66"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -1.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.377"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  +0.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "3.099.70"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  -2.54%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -0.03%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value2 = "Fetch.AI"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "2.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +1.36%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value2 = "ThetaToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.79"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  -4.06%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  +0.65%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  -1.42%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  +0.56%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -2.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "8.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  -0.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "136.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -0.93%  "
